$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 457, pushing the existing rows 457..523
# (and everything below) down by one row (to 458..524).
$ws.Rows("457:457").Insert()

# Populate the newly inserted row 457 with the new data record.
$ws.Cells.Item(457, 1).Value = 10
$ws.Cells.Item(457, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(457, 3).Value = "La Araucanía"
$ws.Cells.Item(457, 4).Value = 44984
$ws.Cells.Item(457, 5).Value = 9
$ws.Cells.Item(457, 6).Value = 100112040
$ws.Cells.Item(457, 7).Value = "Cilantro"
$ws.Cells.Item(457, 8).Value = "Sin especificar"
$ws.Cells.Item(457, 9).Value = "Primera"
$ws.Cells.Item(457, 10).Value = 140
$ws.Cells.Item(457, 11).Value = 7000
$ws.Cells.Item(457, 12).Value = 8000
$ws.Cells.Item(457, 13).Value = 7571
$ws.Cells.Item(457, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(457, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(457, 16).Value = 3786
$ws.Cells.Item(457, 17).Value = 2
$ws.Cells.Item(457, 18).Value = "Hortaliza"
